$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "FuelTanks"

# Update cell values (replace formula in A2 with a plain value)
$ws.Range("A1").Value = 10
$ws.Range("B1").Value = 2
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 2

# Update the selection / active cell
$ws.Range("A3").Select()
